# Updates leve-profit calculations (H/I/J/K/L/M/N columns) across several
# sheets to reflect refreshed market-board pricing data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 2082.7368
$ws.Range("J18").Value = 1122.25
$ws.Range("L18").Value = 1122.25
$ws.Range("N18").Value = -1690.25
$ws.Range("H64").Value = 24305.3
$ws.Range("I64").Value = 47589.777
$ws.Range("K64").Value = 47589.777
$ws.Range("M64").Value = -47341.777
$ws.Range("H67").Value = 24305.3
$ws.Range("I67").Value = 47589.777
$ws.Range("K67").Value = 47589.777
$ws.Range("M67").Value = -46731.777
$ws.Range("H76").Value = 5750
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685
$ws.Range("H79").Value = 5750
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908
$ws.Range("H125").Value = 3300.2
$ws.Range("I125").Value = 768
$ws.Range("J125").Value = 4385.4287
$ws.Range("K125").Value = 6912
$ws.Range("L125").Value = 39468.85830000001
$ws.Range("M125").Value = -4452
$ws.Range("N125").Value = -44388.85830000001
$ws.Range("H135").Value = 3714.6667
$ws.Range("J135").Value = 2249.75
$ws.Range("L135").Value = 20247.75
$ws.Range("N135").Value = -25317.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2686.1365
$ws.Range("I32").Value = 2617
$ws.Range("K32").Value = 2617
$ws.Range("M32").Value = -2330
$ws.Range("H61").Value = 8249.727999999999
$ws.Range("I61").Value = 10801.923
$ws.Range("J61").Value = 4563.222
$ws.Range("K61").Value = 10801.923
$ws.Range("L61").Value = 4563.222
$ws.Range("M61").Value = -10589.923
$ws.Range("N61").Value = -4987.222
$ws.Range("H74").Value = 121498
$ws.Range("I74").Value = 132088.73
$ws.Range("K74").Value = 132088.73
$ws.Range("M74").Value = -131214.73
$ws.Range("H77").Value = 121498
$ws.Range("I77").Value = 132088.73
$ws.Range("K77").Value = 660443.65
$ws.Range("M77").Value = -656075.65
$ws.Range("H102").Value = 4588.151
$ws.Range("I102").Value = 2788.8484
$ws.Range("J102").Value = 7557
$ws.Range("K102").Value = 2788.8484
$ws.Range("L102").Value = 7557
$ws.Range("M102").Value = -1166.8484
$ws.Range("N102").Value = -10801
$ws.Range("H132").Value = 2726.0386
$ws.Range("I132").Value = 2318.907
$ws.Range("K132").Value = 6956.721
$ws.Range("M132").Value = -4426.721
$ws.Range("H136").Value = 8249.727999999999
$ws.Range("I136").Value = 10801.923
$ws.Range("J136").Value = 4563.222
$ws.Range("K136").Value = 32405.769
$ws.Range("L136").Value = 13689.666
$ws.Range("M136").Value = -29855.769
$ws.Range("N136").Value = -18789.666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 9681.706
$ws.Range("I105").Value = 11383
$ws.Range("K105").Value = 11383
$ws.Range("M105").Value = -9636
$ws.Range("H134").Value = 12339.533
$ws.Range("I134").Value = 13985.75
$ws.Range("K134").Value = 41957.25
$ws.Range("M134").Value = -39422.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 559547
$ws.Range("I99").Value = 628240.4
$ws.Range("K99").Value = 628240.4
$ws.Range("M99").Value = -626742.4
$ws.Range("H126").Value = 559547
$ws.Range("I126").Value = 628240.4
$ws.Range("K126").Value = 1884721.2
$ws.Range("M126").Value = -1882251.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H28").Value = 1500
$ws.Range("I28").Value = 1500
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 4500
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -4268
$ws.Range("N28").ClearContents()
$ws.Range("H33").Value = 258.6
$ws.Range("J33").Value = 317.36365
$ws.Range("L33").Value = 1904.1819
$ws.Range("N33").Value = -2470.1819
$ws.Range("H59").Value = 3648.1667
$ws.Range("I59").Value = 2304.2
$ws.Range("K59").Value = 6912.599999999999
$ws.Range("M59").Value = -6372.599999999999
$ws.Range("H86").Value = 2129.2222
$ws.Range("J86").Value = 2420.818
$ws.Range("L86").Value = 7262.454000000001
$ws.Range("N86").Value = -9634.454000000002
$ws.Range("H89").Value = 2129.2222
$ws.Range("J89").Value = 2420.818
$ws.Range("L89").Value = 21787.362
$ws.Range("N89").Value = -33643.362
$ws.Range("H97").Value = 38680.375
$ws.Range("I97").Value = 120890
$ws.Range("J97").Value = 1312.3636
$ws.Range("K97").Value = 362670
$ws.Range("L97").Value = 3937.0908
$ws.Range("M97").Value = -362174
$ws.Range("N97").Value = -4929.0908
$ws.Range("H107").Value = 2214.5454
$ws.Range("J107").Value = 3194.2856
$ws.Range("L107").Value = 9582.856800000001
$ws.Range("N107").Value = -13422.8568
$ws.Range("H131").Value = 71435120
$ws.Range("I131").Value = 166679730
$ws.Range("J131").Value = 1672
$ws.Range("K131").Value = 500039190
$ws.Range("L131").Value = 5016
$ws.Range("M131").Value = -500034150
$ws.Range("N131").Value = -15096
$ws.Range("H137").Value = 7291.967
$ws.Range("I137").Value = 2648.6667
$ws.Range("J137").Value = 10387.5
$ws.Range("K137").Value = 7946.000100000001
$ws.Range("L137").Value = 31162.5
$ws.Range("M137").Value = -2846.000100000001
$ws.Range("N137").Value = -41362.5
$ws.Range("H139").Value = 3336702
$ws.Range("I139").Value = 6001904
$ws.Range("K139").Value = 18005712
$ws.Range("M139").Value = -18000572

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7754.1035
$ws.Range("I70").Value = 7878.45
$ws.Range("K70").Value = 7878.45
$ws.Range("M70").Value = -7608.45
$ws.Range("H73").Value = 7754.1035
$ws.Range("I73").Value = 7878.45
$ws.Range("K73").Value = 7878.45
$ws.Range("M73").Value = -6942.45
$ws.Range("H80").Value = 4415
$ws.Range("I80").Value = 4317.8887
$ws.Range("J80").Value = 4560.6665
$ws.Range("K80").Value = 4317.8887
$ws.Range("L80").Value = 4560.6665
$ws.Range("M80").Value = -3319.8887
$ws.Range("N80").Value = -6556.6665
$ws.Range("H83").Value = 4415
$ws.Range("I83").Value = 4317.8887
$ws.Range("J83").Value = 4560.6665
$ws.Range("K83").Value = 21589.4435
$ws.Range("L83").Value = 22803.3325
$ws.Range("M83").Value = -16597.4435
$ws.Range("N83").Value = -32787.3325
$ws.Range("H122").Value = 12072.637
$ws.Range("I122").Value = 12266.619
$ws.Range("K122").Value = 36799.857
$ws.Range("M122").Value = -34349.857
$ws.Range("H132").Value = 3608.4412
$ws.Range("I132").Value = 3172.3635
$ws.Range("K132").Value = 9517.0905
$ws.Range("M132").Value = -6987.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 63570.43
$ws.Range("I7").Value = 85098.8
$ws.Range("K7").Value = 85098.8
$ws.Range("M7").Value = -84986.8
$ws.Range("H40").Value = 43158.25
$ws.Range("I40").Value = 49010.938
$ws.Range("K40").Value = 49010.938
$ws.Range("M40").Value = -48874.938
$ws.Range("H82").Value = 2148.5173
$ws.Range("I82").Value = 2729.1538
$ws.Range("J82").Value = 1676.75
$ws.Range("K82").Value = 2729.1538
$ws.Range("L82").Value = 1676.75
$ws.Range("M82").Value = -2368.1538
$ws.Range("N82").Value = -2398.75
$ws.Range("H85").Value = 2148.5173
$ws.Range("I85").Value = 2729.1538
$ws.Range("J85").Value = 1676.75
$ws.Range("K85").Value = 2729.1538
$ws.Range("L85").Value = 1676.75
$ws.Range("M85").Value = -1481.1538
$ws.Range("N85").Value = -4172.75
$ws.Range("H126").Value = 63570.43
$ws.Range("I126").Value = 85098.8
$ws.Range("K126").Value = 255296.4
$ws.Range("M126").Value = -252826.4

